$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Price (column D) values are text in this sheet (e.g. "24.213.91" is not a
# valid number). A leading apostrophe forces Excel to keep the literal text
# instead of silently parsing/reformatting it as a number.

# --- Rows 2-49: refreshed Price / Volume(1h) figures ---
$ws.Range("D2").Value = "'24.200.63"
$ws.Range("E2").Value = "  -1.82%  "
$ws.Range("D3").Value = "'1.648.83"
$ws.Range("E3").Value = "  -1.74%  "
$ws.Range("D4").Value = "'0.9998"
$ws.Range("E4").Value = "  -0.77%  "
$ws.Range("D5").Value = "'308.78"
$ws.Range("E5").Value = "  -1.42%  "
$ws.Range("D6").Value = "'1.000"
$ws.Range("E6").Value = "  -0.45%  "
$ws.Range("D7").Value = "'0.3929"
$ws.Range("E7").Value = "  +0.25%  "
$ws.Range("D8").Value = "'0.3869"
$ws.Range("E8").Value = "  -2.34%  "
$ws.Range("D9").Value = "'1.001"
$ws.Range("E9").Value = "  -0.65%  "
$ws.Range("D10").Value = "'1.368"
$ws.Range("E10").Value = "  -3.55%  "
$ws.Range("D11").Value = "'49.51"
$ws.Range("E11").Value = "  -4.58%  "
$ws.Range("D12").Value = "'0.08611"
$ws.Range("E12").Value = "  -0.72%  "
$ws.Range("D13").Value = "'23.67"
$ws.Range("E13").Value = "  -6.14%  "
$ws.Range("D14").Value = "'7.120"
$ws.Range("E14").Value = "  -2.76%  "
$ws.Range("D15").Value = "'0.00001294"
$ws.Range("E15").Value = "  -1.87%  "
$ws.Range("D16").Value = "'7.509"
$ws.Range("E16").Value = "  -3.74%  "
$ws.Range("D17").Value = "'1.650.81"
$ws.Range("E17").Value = "  -1.35%  "
$ws.Range("D18").Value = "'94.94"
$ws.Range("E18").Value = "  +1.50%  "
$ws.Range("D20").Value = "'20.42"
$ws.Range("E20").Value = "  +0.86%  "
$ws.Range("D21").Value = "'6.917"
$ws.Range("E21").Value = "  -2.96%  "
$ws.Range("D22").Value = "'1.000"
$ws.Range("E22").Value = "  -0.60%  "
$ws.Range("D23").Value = "'13.61"
$ws.Range("E23").Value = "  -2.89%  "
$ws.Range("D24").Value = "'24.194.95"
$ws.Range("E24").Value = "  -1.86%  "
$ws.Range("D25").Value = "'2.439"
$ws.Range("E25").Value = "  +3.79%  "
$ws.Range("D26").Value = "'2.867"
$ws.Range("E26").Value = "  +3.78%  "
$ws.Range("D27").Value = "'22.44"
$ws.Range("E27").Value = "  -4.03%  "
$ws.Range("D28").Value = "'158.46"
$ws.Range("E28").Value = "  -2.09%  "
$ws.Range("D29").Value = "'8.403"
$ws.Range("E29").Value = "  +7.23%  "
$ws.Range("D30").Value = "'5.398"
$ws.Range("E30").Value = "  -5.32%  "
$ws.Range("D31").Value = "'141.01"
$ws.Range("E31").Value = "  -5.41%  "
$ws.Range("D32").Value = "'2.418"
$ws.Range("E32").Value = "  -4.70%  "
$ws.Range("D33").Value = "'1.832.13"
$ws.Range("E33").Value = "  -1.27%  "
$ws.Range("D36").Value = "'0.02919"
$ws.Range("E36").Value = "  -5.09%  "
$ws.Range("D37").Value = "'0.2702"
$ws.Range("E37").Value = "  -3.26%  "
$ws.Range("D38").Value = "'0.9591"
$ws.Range("E38").Value = "  -4.94%  "
$ws.Range("D39").Value = "'0.09227"
$ws.Range("E39").Value = "  -3.21%  "
$ws.Range("D40").Value = "'10.41"
$ws.Range("E40").Value = "  -1.31%  "
$ws.Range("D41").Value = "'1.460"
$ws.Range("E41").Value = "  -1.09%  "
$ws.Range("D42").Value = "'0.7560"
$ws.Range("E42").Value = "  -4.50%  "
$ws.Range("D43").Value = "'13.06"
$ws.Range("E43").Value = "  -4.26%  "
$ws.Range("D44").Value = "'16.21"
$ws.Range("E44").Value = "  -2.86%  "
$ws.Range("D45").Value = "'0.6923"
$ws.Range("E45").Value = "  -2.99%  "
$ws.Range("D46").Value = "'2.470"
$ws.Range("E46").Value = "  -4.06%  "
$ws.Range("D47").Value = "'4.097"
$ws.Range("E47").Value = "  -1.98%  "
$ws.Range("D48").Value = "'0.9997"
$ws.Range("E48").Value = "  -0.43%  "
$ws.Range("D49").Value = "'0.08391"
$ws.Range("E49").Value = "  -3.12%  "

# --- Row 19: Volume(1h) only (Price unchanged) ---
$ws.Range("E19").Value = "  -2.89%  "

# --- Rows 34/35: InternetComputer(DFINITY) and Hedera swap rank order ---
$ws.Range("B34").Value = "Hedera"
$ws.Range("C34").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D34").Value = "'0.08159"
$ws.Range("E34").Value = "  -3.14%  "

$ws.Range("B35").Value = "InternetComputer(DFINITY)"
$ws.Range("C35").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D35").Value = "'6.949"
$ws.Range("E35").Value = "  +0.31%  "

# --- Rows 50/51: Quant and Flow swap rank order ---
$ws.Range("B50").Value = "Flow"
$ws.Range("C50").Value = "https://coinranking.com/coin/QQ0NCmjVq+flow-flow"
$ws.Range("D50").Value = "'1.265"
$ws.Range("E50").Value = "  -5.15%  "

$ws.Range("B51").Value = "Quant"
$ws.Range("C51").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D51").Value = "'133.87"
$ws.Range("E51").Value = "  -3.25%  "

